# Apply "Results - second part" edits to the Test sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing II deo scores for several rows (formulas in column G
# recalculate automatically since they already reference C:F).
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 26
$ws.Range("F5").Value = 19

$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 28
$ws.Range("F9").Value = 19

$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = 20

$ws.Range("D13").Value = 18
$ws.Range("E13").Value = 26
$ws.Range("F13").Value = 14

$ws.Range("D14").Value = 18
$ws.Range("E14").Value = 24
$ws.Range("F14").Value = 19

$ws.Range("D16").Value = 20
$ws.Range("E16").Value = 26
$ws.Range("F16").Value = 20

# Move the active selection to A11 as recorded in the saved view state.
$ws.Range("A11").Select()
